$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the generation Date value ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-06-30T08:07:26+01:00"

# --- Elements sheet: fill in the Binding Strength / Description / Value Set
#     for the Extension.value[x] row (row 6) ---
$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("X6").Value = "required"
$wsElem.Range("Y6").Value = ""
$wsElem.Range("Z6").Value = "http://nphcda.gov.ng/ig/fhir/ValueSet/ng-sibling-health-status"

# Column Z (Binding Value Set) now holds a much longer value, so re-fit its width
$wsElem.Columns.Item(26).ColumnWidth = 50.014322916666664
